$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("emp")

# --- Fix D83 / D84: "NA" -> "N/A" (quote-prefixed so Excel keeps it as text,
#     mirroring the existing quotePrefix formatting already used on those cells) ---
$ws.Range("D83").Value = "'N/A"
$ws.Range("D84").Value = "'N/A"

# --- Append the 31 new rows (99-129) to the Tabla3 table ---
$lo = $ws.ListObjects.Item("Tabla3")

$newRows = @(
    @{B="CO4828076144X"; G=100; H=1;  J=79;   K=107},
    @{B="CO4828076145X"; G=100; H=1;  J=83;   K=111},
    @{B="CO4840076160X"; G=50;  H=2;  J=86;   K=114},
    @{B="CO4806076200X"; G=78;  H=2;  J=278;  K=306},
    @{B="CO483T080160X"; G=84;  H=9;  J=1856; K=1884},
    @{B="CO483U080160X"; G=84;  H=9;  J=1936; K=1964},
    @{B="CO483V080254X"; G=48;  H=14; J=1848; K=1876},
    @{B="CO484F075350X"; G=39;  H=10; J=1348; K=1376},
    @{B="CO481X767160X"; G=91;  H=32; J=1860; K=1888},
    @{B="CO481Z767160X"; G=91;  H=17; J=1904; K=1932},
    @{B="CO483S080160X"; G=84;  H=11; J=1760; K=1788},
    @{B="CO481Y767160X"; G=91;  H=25; J=1886; K=1914},
    @{B="CO485B762170X"; G=91;  H=10; J=1850; K=1878},
    @{B="CO484M638127X"; G=135; H=20; J=1890; K=1918},
    @{B="CO4807076160X"; G=91;  H=16; J=1785; K=1813},
    @{B="CO481L076150X"; G=104; H=19; J=1800; K=1828},
    @{B="CO481N767150X"; G=104; H=15; J=1792; K=1820},
    @{B="CO481L767160X"; G=91;  H=19; J=1728; K=1756},
    @{B="CO485Q060140X"; G=128; H=31; J=1880; K=1908},
    @{B="CO484V603085X"; G=210; H=19; J=1720; K=1748},
    @{B="CO487F628127";  G=153; H=13; J=1860; K=1888},
    @{B="CO481K076150X"; G=98;  H=30; J=1770; K=1798},
    @{B="CO487N774188X"; G=84;  H=16; J=1937; K=1965},
    @{B="CO489F768152X"; G=98;  H=28; J=1848; K=1876},
    @{B="CO487W762140X"; G=112; H=16; J=1868; K=1896},
    @{B="CO4807076144X"; G=100; H=1;  J=103;  K=131},
    @{B="CO481K076150X"; G=100; H=1;  J=59;   K=87},
    @{B="CO4824076143X"; G=100; H=1;  J=114;  K=142},
    @{B="CO4841076180X"; G=100; H=1;  J=107;  K=135},
    @{B="CO4843076144X"; G=100; H=1;  J=107;  K=135},
    @{B="CO4824076144X"; G=100; H=1;  J=92;   K=120}
)

foreach ($data in $newRows) {
    $row = $lo.ListRows.Add()
    $r = $lo.Range.Rows.Count + $lo.Range.Row - 1

    $ws.Range("A$r").Value = "TACON"
    $ws.Range("B$r").Value = $data.B
    $ws.Range("C$r").Value = "RECICLADA"
    $ws.Range("D$r").Value = "'N/A"
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("I$r").Formula = "=Tabla3[[#This Row],[MxC]]*Tabla3[[#This Row],[CAMAS]]"
    $ws.Range("J$r").Value = $data.J
    $ws.Range("K$r").Value = $data.K
    $ws.Range("L$r").Value = "LINK"
    $ws.Range("M$r").Value = $false
}

$ws.Range("E99").Select()
